$d = $word.ActiveDocument

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pBdr><w:top w:val="single" w:sz="4" w:space="1" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="4" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="1" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="4" w:color="auto"/></w:pBdr><w:tabs><w:tab w:val="left" w:leader="underscore" w:pos="10800"/></w:tabs><w:spacing w:before="240"/><w:contextualSpacing/></w:pPr><w:r><w:t xml:space="preserve">Member’s Name: </w:t></w:r><w:r><w:tab/></w:r></w:p><w:p><w:pPr><w:pBdr><w:top w:val="single" w:sz="4" w:space="1" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="4" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="1" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="4" w:color="auto"/></w:pBdr><w:tabs><w:tab w:val="left" w:leader="underscore" w:pos="10800"/></w:tabs><w:contextualSpacing/></w:pPr><w:r><w:t xml:space="preserve">4-H Club Name: </w:t></w:r><w:r><w:tab/></w:r></w:p><w:p><w:pPr><w:pBdr><w:top w:val="single" w:sz="4" w:space="1" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="4" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="1" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="4" w:color="auto"/></w:pBdr><w:tabs><w:tab w:val="left" w:leader="underscore" w:pos="10800"/></w:tabs><w:contextualSpacing/></w:pPr><w:r><w:t xml:space="preserve">Name of Club Leader: </w:t></w:r><w:r><w:tab/></w:r></w:p><w:p><w:pPr><w:pBdr><w:top w:val="single" w:sz="4" w:space="1" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="4" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="1" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="4" w:color="auto"/></w:pBdr><w:tabs><w:tab w:val="left" w:leader="underscore" w:pos="3240"/><w:tab w:val="left" w:leader="underscore" w:pos="6570"/><w:tab w:val="left" w:leader="underscore" w:pos="9360"/></w:tabs><w:contextualSpacing/></w:pPr><w:r><w:t xml:space="preserve">Report Year: </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve">Year(s) Enrolled in 4-H: </w:t></w:r><w:r><w:tab/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$r = $d.Range(0, 0)
$r.InsertXML($xml)

Write-Output ("ParaCount=" + $d.Paragraphs.Count)
